$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Snapshot current (pre-edit) row values for columns D,M,N,O,P,Q,R,S,T
# before any writes happen, since the permutation reads from multiple source rows.
$src = @{}
$src[2] = @{ "D" = 44616; "M" = 70; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos empedrada"; "R" = "Provincia de Limarí" }
$src[3] = @{ "D" = 44630; "M" = 75; "N" = 15000; "O" = 15000; "P" = 15000; "S" = 1071; "T" = 14; "Q" = "`$/caja 14 kilos empedrada"; "R" = "Provincia de Limarí" }
$src[4] = @{ "D" = 44585; "M" = 50; "N" = 22500; "O" = 22500; "P" = 22500; "S" = 1500; "T" = 15; "Q" = "`$/caja 15 kilos empedrada"; "R" = "Provincia de Limarí" }
$src[5] = @{ "D" = 44239; "M" = 70; "N" = 15000; "O" = 15000; "P" = 15000; "S" = 1000; "T" = 15; "Q" = "`$/caja 15 kilos granel"; "R" = "Provincia de Limarí" }
$src[6] = @{ "D" = 44259; "M" = 80; "N" = 12000; "O" = 12000; "P" = 12000; "S" = 800; "T" = 15; "Q" = "`$/caja 15 kilos empedrada"; "R" = "Provincia de Limarí" }
$src[7] = @{ "D" = 44270; "M" = 85; "N" = 12000; "O" = 12000; "P" = 12000; "S" = 857; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia del Elquí" }
$src[8] = @{ "D" = 44614; "M" = 54; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }
$src[9] = @{ "D" = 44278; "M" = 45; "N" = 13000; "O" = 13000; "P" = 13000; "S" = 929; "T" = 14; "Q" = "`$/caja 14 kilos empedrada"; "R" = "Provincia del Elquí" }
$src[10] = @{ "D" = 44314; "M" = 56; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }
$src[11] = @{ "D" = 44260; "M" = 56; "N" = 13000; "O" = 13000; "P" = 13000; "S" = 929; "T" = 14; "Q" = "`$/caja 14 kilos empedrada"; "R" = "Provincia del Elquí" }
$src[12] = @{ "D" = 44245; "M" = 50; "N" = 15000; "O" = 15000; "P" = 15000; "S" = 1000; "T" = 15; "Q" = "`$/caja 15 kilos granel"; "R" = "Provincia de Limarí" }
$src[13] = @{ "D" = 44323; "M" = 60; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }
$src[14] = @{ "D" = 44592; "M" = 54; "N" = 20000; "O" = 20000; "P" = 20000; "S" = 1333; "T" = 15; "Q" = "`$/caja 15 kilos empedrada"; "R" = "Provincia de Limarí" }
$src[15] = @{ "D" = 44316; "M" = 48; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }
$src[16] = @{ "D" = 44322; "M" = 50; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }
$src[17] = @{ "D" = 44588; "M" = 85; "N" = 19000; "O" = 20000; "P" = 19529; "S" = 1395; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }
$src[18] = @{ "D" = 44320; "M" = 45; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }
$src[19] = @{ "D" = 44313; "M" = 36; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }
$src[20] = @{ "D" = 44238; "M" = 60; "N" = 15000; "O" = 15000; "P" = 15000; "S" = 1000; "T" = 15; "Q" = "`$/caja 15 kilos granel"; "R" = "Provincia de Limarí" }
$src[21] = @{ "D" = 44242; "M" = 45; "N" = 12000; "O" = 12000; "P" = 12000; "S" = 800; "T" = 15; "Q" = "`$/caja 15 kilos granel"; "R" = "Provincia de Limarí" }
$src[22] = @{ "D" = 44252; "M" = 60; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos empedrada"; "R" = "Provincia de Limarí" }
$src[23] = @{ "D" = 44271; "M" = 50; "N" = 12000; "O" = 12000; "P" = 12000; "S" = 857; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia del Elquí" }
$src[24] = @{ "D" = 44315; "M" = 65; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }
$src[25] = @{ "D" = 44627; "M" = 56; "N" = 17000; "O" = 17000; "P" = 17000; "S" = 1214; "T" = 14; "Q" = "`$/caja 14 kilos empedrada"; "R" = "Provincia de Limarí" }
$src[26] = @{ "D" = 44312; "M" = 68; "N" = 14000; "O" = 14000; "P" = 14000; "S" = 1000; "T" = 14; "Q" = "`$/caja 14 kilos granel"; "R" = "Provincia de Limarí" }

# Apply permutation: row r gets the D/M/N/O/P/Q/R/S/T tuple that previously lived on row mapping[r]
$mapping = @{}
$mapping[2] = 17
$mapping[3] = 24
$mapping[4] = 10
$mapping[5] = 20
$mapping[6] = 16
$mapping[7] = 3
$mapping[8] = 23
$mapping[9] = 18
$mapping[10] = 21
$mapping[11] = 6
$mapping[12] = 11
$mapping[13] = 2
$mapping[14] = 15
$mapping[15] = 25
$mapping[16] = 5
$mapping[17] = 8
$mapping[18] = 22
$mapping[19] = 12
$mapping[20] = 9
$mapping[21] = 4
$mapping[22] = 19
$mapping[23] = 13
$mapping[24] = 7
$mapping[25] = 26
$mapping[26] = 14

foreach ($r in ($mapping.Keys | Sort-Object)) {
    $s = $src[$mapping[$r]]
    $ws.Cells.Item($r, 4).Value = $s["D"]
    $ws.Cells.Item($r, 13).Value = $s["M"]
    $ws.Cells.Item($r, 14).Value = $s["N"]
    $ws.Cells.Item($r, 15).Value = $s["O"]
    $ws.Cells.Item($r, 16).Value = $s["P"]
    $ws.Cells.Item($r, 17).Value = $s["Q"]
    $ws.Cells.Item($r, 18).Value = $s["R"]
    $ws.Cells.Item($r, 19).Value = $s["S"]
    $ws.Cells.Item($r, 20).Value = $s["T"]
}

Write-Host "Done applying permutation to rows 2-26"
